$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.417.73"
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = "'1.643.46"
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = "'0.9986"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").Value = "'299.40"
$ws.Range("E6").Value = '  -1.21%  '
$ws.Range("D7").Value = "'0.3793"
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("D8").Value = "'50.30"
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").Value = "'0.3489"
$ws.Range("E9").Value = '  -3.39%  '
$ws.Range("D10").Value = "'0.08068"
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("D11").Value = "'1.218"
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = "'0.9987"
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").Value = "'22.06"
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("D14").Value = "'6.308"
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = "'7.267"
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = "'0.00001208"
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = "'1.636.21"
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = "'94.90"
$ws.Range("E18").Value = '  -2.87%  '
$ws.Range("D19").Value = "'0.06963"
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").Value = "'6.621"
$ws.Range("E20").Value = '  -2.49%  '
$ws.Range("D21").Value = "'17.34"
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = "'12.42"
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").Value = "'23.440.56"
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").Value = "'2.433"
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("D26").Value = "'2.980"
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").Value = "'20.99"
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("D28").Value = "'149.66"
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").Value = "'5.174"
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("D30").Value = "'131.62"
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = "'1.820.10"
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").Value = "'6.821"
$ws.Range("E32").Value = '  -4.04%  '
$ws.Range("D33").Value = "'2.134"
$ws.Range("E33").Value = '  -5.33%  '
$ws.Range("D34").Value = "'11.25"
$ws.Range("E34").Value = '  -6.23%  '
$ws.Range("D35").Value = "'0.9911"
$ws.Range("E35").Value = '  -6.15%  '
$ws.Range("D36").Value = "'0.02686"
$ws.Range("E36").Value = '  -4.37%  '
$ws.Range("D37").Value = "'0.08768"
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").Value = "'0.2425"
$ws.Range("E38").Value = '  -3.22%  '
$ws.Range("D39").Value = "'5.884"
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").Value = "'0.06835"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("D41").Value = "'12.79"
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("D42").Value = "'0.6825"
$ws.Range("E42").Value = '  -2.40%  '
$ws.Range("D43").Value = "'1.288"
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("D44").Value = "'15.44"
$ws.Range("E44").Value = '  -3.09%  '
$ws.Range("D45").Value = "'0.9980"
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("D46").Value = "'0.6335"
$ws.Range("E46").Value = '  -2.71%  '
$ws.Range("D47").Value = "'2.239"
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").Value = "'3.907"
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("D49").Value = "'127.05"
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("D50").Value = "'0.07674"
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = '  +2.31%  '
